# Updated faculty data for 3-2 semester with real faculty list
# Fills in the previously-blank "19-Dec" attendance column (M) for
# students in rows 2-11 of the Attendance sheet, and refreshes the
# dependent Total Present / Total Absent / Percentage summary cells
# (DW/DX/DY) to match the new attendance counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

function Set-TextValue($range, $text) {
    # Force the value to be stored as a literal string (avoids Excel
    # auto-converting things like "80.0%" into a percentage number),
    # then clear the formatting so no stray style index gets attached
    # to the cell (keeps the cell looking like the untouched ones).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
$ws.Range("M2").Value = "A"
$ws.Range("DX2").Value = 2
Set-TextValue $ws.Range("DY2") "80.0%"

# Row 3
$ws.Range("E3").Value = "A"
$ws.Range("M3").Value = "P"
$ws.Range("DX3").Value = 2
Set-TextValue $ws.Range("DY3") "80.0%"

# Row 4
$ws.Range("E4").Value = "A"
$ws.Range("M4").Value = "P"
$ws.Range("DX4").Value = 2
Set-TextValue $ws.Range("DY4") "80.0%"

# Row 5
$ws.Range("M5").Value = "A"
$ws.Range("DX5").Value = 2
Set-TextValue $ws.Range("DY5") "80.0%"

# Row 6
$ws.Range("M6").Value = "A"
$ws.Range("DX6").Value = 2
Set-TextValue $ws.Range("DY6") "80.0%"

# Row 7
$ws.Range("M7").Value = "P"
$ws.Range("DW7").Value = 9
Set-TextValue $ws.Range("DY7") "90.0%"

# Row 8
$ws.Range("M8").Value = "P"
$ws.Range("DW8").Value = 9
Set-TextValue $ws.Range("DY8") "90.0%"

# Row 9
$ws.Range("M9").Value = "A"
$ws.Range("DX9").Value = 2
Set-TextValue $ws.Range("DY9") "80.0%"

# Row 10
$ws.Range("M10").Value = "P"
$ws.Range("DW10").Value = 9
Set-TextValue $ws.Range("DY10") "90.0%"

# Row 11
$ws.Range("M11").Value = "P"
$ws.Range("DW11").Value = 9
Set-TextValue $ws.Range("DY11") "90.0%"
